$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update column B values to their simplified (shortened) scores ---
$ws.Range("B2").Value = 199
$ws.Range("B3").Value = 60
$ws.Range("B4").Value = 88
$ws.Range("B5").Value = 481
$ws.Range("B6").Value = 597
$ws.Range("B7").Value = 410
$ws.Range("B8").Value = 880
$ws.Range("B9").Value = 427
$ws.Range("B10").Value = 828
$ws.Range("B11").Value = 391
$ws.Range("B12").Value = 330
$ws.Range("B13").Value = 399
$ws.Range("B14").Value = 989
$ws.Range("B15").Value = 280
$ws.Range("B16").Value = 633
$ws.Range("B17").Value = 807
$ws.Range("B18").Value = 439
$ws.Range("B19").Value = 758
$ws.Range("B20").Value = 3
$ws.Range("B21").Value = 607
$ws.Range("B22").Value = 376
$ws.Range("B23").Value = 674
$ws.Range("B24").Value = 261
$ws.Range("B25").Value = 375
$ws.Range("B26").Value = 933
$ws.Range("B27").Value = 527
$ws.Range("B28").Value = 39
$ws.Range("B29").Value = 214
$ws.Range("B30").Value = 482
$ws.Range("B31").Value = 907
$ws.Range("B32").Value = 543
$ws.Range("B33").Value = 64
$ws.Range("B34").Value = 222
$ws.Range("B35").Value = 903
$ws.Range("B36").Value = 734
$ws.Range("B37").Value = 43
$ws.Range("B38").Value = 445
$ws.Range("B39").Value = 530
$ws.Range("B40").Value = 169
$ws.Range("B41").Value = 446
$ws.Range("B42").Value = 649
$ws.Range("B43").Value = 70
$ws.Range("B44").Value = 482
$ws.Range("B45").Value = 595
$ws.Range("B46").Value = 233
$ws.Range("B47").Value = 721
$ws.Range("B48").Value = 917
$ws.Range("B49").Value = 694
$ws.Range("B50").Value = 269
$ws.Range("B51").Value = 37
$ws.Range("B52").Value = 555
$ws.Range("B53").Value = 966
$ws.Range("B54").Value = 839
$ws.Range("B55").Value = 481
$ws.Range("B56").Value = 908
$ws.Range("B57").Value = 68
$ws.Range("B58").Value = 749
$ws.Range("B59").Value = 473
$ws.Range("B60").Value = 31
$ws.Range("B61").Value = 686
$ws.Range("B62").Value = 70
$ws.Range("B63").Value = 23
$ws.Range("B64").Value = 798
$ws.Range("B65").Value = 974
$ws.Range("B66").Value = 598
$ws.Range("B67").Value = 948
$ws.Range("B68").Value = 516
$ws.Range("B69").Value = 418
$ws.Range("B70").Value = 284
$ws.Range("B71").Value = 983
$ws.Range("B72").Value = 131
$ws.Range("B73").Value = 961
$ws.Range("B74").Value = 882
$ws.Range("B75").Value = 666
$ws.Range("B76").Value = 44
$ws.Range("B77").Value = 64
$ws.Range("B78").Value = 13
$ws.Range("B79").Value = 464
$ws.Range("B80").Value = 825

# --- Narrow column A slightly and widen column B so the shorter scores are
#     still easy to compare visually against column A ---
$ws.Columns.Item(1).ColumnWidth = 13
$ws.Columns.Item(2).ColumnWidth = 23.5

# --- Scroll back up (closer to the top of the data) and select B2 ---
$excel.ActiveWindow.ScrollRow = 31
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B2").Select() | Out-Null
